$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.861.37'
$ws.Range("E2").Value = '  +0.35%  '

$ws.Range("D3").Value = '2.288.37'
$ws.Range("E3").Value = '  -1.64%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '96.97'
$ws.Range("E5").Value = '  +1.05%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '269.62'
$ws.Range("E6").Value = '  -0.80%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.625'
$ws.Range("E7").Value = '  -0.21%  '

$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.609'
$ws.Range("E9").Value = '  -1.89%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.30'
$ws.Range("E10").Value = '  +0.90%  '

$ws.Range("E11").Value = '  -0.56%  '

$ws.Range("E12").Value = '  -2.46%  '

$ws.Range("E13").Value = '  +1.64%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.88'
$ws.Range("E14").Value = '  +1.55%  '

$ws.Range("D15").Value = '2.630.53'
$ws.Range("E15").Value = '  -1.67%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.859'
$ws.Range("E16").Value = '  +0.40%  '

$ws.Range("D17").Value = '2.286.07'
$ws.Range("E17").Value = '  -1.99%  '

$ws.Range("D18").Value = '43.734.99'
$ws.Range("E18").Value = '  +0.24%  '

$ws.Range("E19").Value = '  +3.46%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.19'
$ws.Range("E20").Value = '  -1.99%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.07'
$ws.Range("E21").Value = '  +0.14%  '

$ws.Range("E22").Value = '  +9.87%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.78'
$ws.Range("E23").Value = '  -2.26%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.14'
$ws.Range("E24").Value = '  -4.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.71'
$ws.Range("E25").Value = '  +6.29%  '

$ws.Range("E26").Value = '  +0.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.36'
$ws.Range("E27").Value = '  +0.68%  '

$ws.Range("E28").Value = '  +1.19%  '

$ws.Range("E29").Value = '  -0.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.61'
$ws.Range("E30").Value = '  -0.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.55'
$ws.Range("E31").Value = '  +1.65%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.78'
$ws.Range("E32").Value = '  -3.10%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0894'
$ws.Range("E33").Value = '  -0.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.43'
$ws.Range("E34").Value = '  -1.02%  '

$ws.Range("E35").Value = '  +0.35%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.76'
$ws.Range("E36").Value = '  +8.94%  '

$ws.Range("E37").Value = '  -0.05%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0352'
$ws.Range("E38").Value = '  -1.84%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.44'
$ws.Range("E39").Value = '  +1.30%  '

$ws.Range("E40").Value = '  +0.64%  '

$ws.Range("E41").Value = '  -2.19%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.27'
$ws.Range("E42").Value = '  +1.70%  '

$ws.Range("E43").Value = '  -0.31%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '64.36'
$ws.Range("E44").Value = '  +4.05%  '

$ws.Range("E45").Value = '  -3.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.71'
$ws.Range("E46").Value = '  -4.17%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '98.34'
$ws.Range("E48").Value = '  -2.39%  '

$ws.Range("E49").Value = '  -1.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.442'
$ws.Range("E50").Value = '  +6.71%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.52'
$ws.Range("E51").Value = '  +11.67%  '
